$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1490.8966
$ws.Range("J17").Value = 1564.7307
$ws.Range("L17").Value = 4694.1921
$ws.Range("N17").Value = -5030.1921
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H80").Value = 1017.2
$ws.Range("I80").Value = 523.3333
$ws.Range("J80").Value = 1758
$ws.Range("K80").Value = 1569.9999
$ws.Range("L80").Value = 5274
$ws.Range("M80").Value = -571.9999
$ws.Range("N80").Value = -7270
$ws.Range("H83").Value = 1017.2
$ws.Range("I83").Value = 523.3333
$ws.Range("J83").Value = 1758
$ws.Range("K83").Value = 4709.9997
$ws.Range("L83").Value = 15822
$ws.Range("M83").Value = 282.0002999999997
$ws.Range("N83").Value = -25806
$ws.Range("H88").Value = 7147133.5
$ws.Range("I88").Value = 16669249
$ws.Range("J88").Value = 5546.625
$ws.Range("K88").Value = 16669249
$ws.Range("L88").Value = 5546.625
$ws.Range("M88").Value = -16668843
$ws.Range("N88").Value = -6358.625
$ws.Range("H91").Value = 7147133.5
$ws.Range("I91").Value = 16669249
$ws.Range("J91").Value = 5546.625
$ws.Range("K91").Value = 16669249
$ws.Range("L91").Value = 5546.625
$ws.Range("M91").Value = -16667845
$ws.Range("N91").Value = -8354.625
$ws.Range("H98").Value = 2264.6333
$ws.Range("I98").Value = 2505.6667
$ws.Range("K98").Value = 2505.6667
$ws.Range("M98").Value = -1007.6667
$ws.Range("H122").Value = 2264.6333
$ws.Range("I122").Value = 2505.6667
$ws.Range("K122").Value = 7517.000100000001
$ws.Range("M122").Value = -5067.000100000001
$ws.Range("H137").Value = 1888.5209
$ws.Range("I137").Value = 1221.0526
$ws.Range("J137").Value = 4424.9
$ws.Range("K137").Value = 3663.1578
$ws.Range("L137").Value = 13274.7
$ws.Range("M137").Value = -1113.1578
$ws.Range("N137").Value = -18374.7
$ws.Range("H138").Value = 3103.7036
$ws.Range("I138").Value = 2795.4
$ws.Range("K138").Value = 8386.200000000001
$ws.Range("M138").Value = -3246.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1610.0597
$ws.Range("I32").Value = 1610.0597
$ws.Range("K32").Value = 1610.0597
$ws.Range("M32").Value = -1323.0597
$ws.Range("H74").Value = 2913.3618
$ws.Range("I74").Value = 2481.4285
$ws.Range("J74").Value = 4173.1665
$ws.Range("K74").Value = 2481.4285
$ws.Range("L74").Value = 4173.1665
$ws.Range("M74").Value = -1607.4285
$ws.Range("N74").Value = -5921.1665
$ws.Range("H77").Value = 2913.3618
$ws.Range("I77").Value = 2481.4285
$ws.Range("J77").Value = 4173.1665
$ws.Range("K77").Value = 12407.1425
$ws.Range("L77").Value = 20865.8325
$ws.Range("M77").Value = -8039.1425
$ws.Range("N77").Value = -29601.8325
$ws.Range("H132").Value = 2642.923
$ws.Range("J132").Value = 3324.75
$ws.Range("L132").Value = 9974.25
$ws.Range("N132").Value = -15034.25
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3833.2
$ws.Range("I86").Value = 3163
$ws.Range("K86").Value = 3163
$ws.Range("M86").Value = -2040
$ws.Range("H89").Value = 3833.2
$ws.Range("I89").Value = 3163
$ws.Range("K89").Value = 15815
$ws.Range("M89").Value = -10199
$ws.Range("H134").Value = 2348.147
$ws.Range("I134").Value = 1804.7693
$ws.Range("K134").Value = 5414.3079
$ws.Range("M134").Value = -2879.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5000.0415
$ws.Range("I31").Value = 4588.125
$ws.Range("J31").Value = 5206
$ws.Range("K31").Value = 4588.125
$ws.Range("L31").Value = 5206
$ws.Range("M31").Value = -4293.125
$ws.Range("N31").Value = -5796
$ws.Range("H34").Value = 5000.0415
$ws.Range("I34").Value = 4588.125
$ws.Range("J34").Value = 5206
$ws.Range("K34").Value = 4588.125
$ws.Range("L34").Value = 5206
$ws.Range("M34").Value = -4386.125
$ws.Range("N34").Value = -5610
$ws.Range("H50").Value = 33518.4
$ws.Range("J50").Value = 33518.4
$ws.Range("L50").Value = 33518.4
$ws.Range("N50").Value = -34768.4
$ws.Range("H51").Value = 59098
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 2203.08
$ws.Range("I58").Value = 784.0833
$ws.Range("J58").Value = 3512.923
$ws.Range("K58").Value = 784.0833
$ws.Range("L58").Value = 3512.923
$ws.Range("M58").Value = -581.0833
$ws.Range("N58").Value = -3918.923
$ws.Range("H59").Value = 65777.5
$ws.Range("J59").Value = 65777.5
$ws.Range("L59").Value = 65777.5
$ws.Range("N59").Value = -68067.5
$ws.Range("H60").Value = 25000
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26022
$ws.Range("H61").Value = 59098
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H132").Value = 18526690
$ws.Range("J132").Value = 47633120
$ws.Range("L132").Value = 142899360
$ws.Range("N132").Value = -142904420
$ws.Range("H134").Value = 2165.8108
$ws.Range("I134").Value = 2014.5938
$ws.Range("J134").Value = 3133.6
$ws.Range("K134").Value = 6043.7814
$ws.Range("L134").Value = 9400.799999999999
$ws.Range("M134").Value = -3508.7814
$ws.Range("N134").Value = -14470.8
$ws.Range("H136").Value = 2203.08
$ws.Range("I136").Value = 784.0833
$ws.Range("J136").Value = 3512.923
$ws.Range("K136").Value = 2352.2499
$ws.Range("L136").Value = 10538.769
$ws.Range("M136").Value = 197.7501000000002
$ws.Range("N136").Value = -15638.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 126289.17
$ws.Range("I4").Value = 109566.164
$ws.Range("J4").Value = 666999.7
$ws.Range("K4").Value = 328698.492
$ws.Range("L4").Value = 2000999.1
$ws.Range("M4").Value = -328586.492
$ws.Range("N4").Value = -2001223.1
$ws.Range("H11").Value = 25502500
$ws.Range("I11").Value = 25502500
$ws.Range("K11").Value = 76507500
$ws.Range("M11").Value = -76507360
$ws.Range("H34").Value = 367.8889
$ws.Range("I34").Value = 288.875
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 866.625
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -782.625
$ws.Range("N34").Value = -3168
$ws.Range("H126").Value = 7000
$ws.Range("I126").Value = 7000
$ws.Range("K126").Value = 21000
$ws.Range("M126").Value = -16060
$ws.Range("H132").Value = 1818.25
$ws.Range("J132").Value = 2003
$ws.Range("L132").Value = 18027
$ws.Range("N132").Value = -23087
$ws.Range("H137").Value = 2776.3076
$ws.Range("I137").Value = 2376.625
$ws.Range("K137").Value = 7129.875
$ws.Range("M137").Value = -2029.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 1500
$ws.Range("J23").Value = 1500
$ws.Range("L23").Value = 1500
$ws.Range("N23").Value = -1946
$ws.Range("H126").Value = 6260.1904
$ws.Range("J126").Value = 9350.416999999999
$ws.Range("L126").Value = 28051.251
$ws.Range("N126").Value = -32991.251
$ws.Range("H132").Value = 2202.5898
$ws.Range("I132").Value = 1800.8462
$ws.Range("J132").Value = 3006.077
$ws.Range("K132").Value = 5402.5386
$ws.Range("L132").Value = 9018.231
$ws.Range("M132").Value = -2872.5386
$ws.Range("N132").Value = -14078.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4869.5386
$ws.Range("I7").Value = 3855.375
$ws.Range("J7").Value = 6492.2
$ws.Range("K7").Value = 3855.375
$ws.Range("L7").Value = 6492.2
$ws.Range("M7").Value = -3743.375
$ws.Range("N7").Value = -6716.2
$ws.Range("H40").Value = 49644.523
$ws.Range("I40").Value = 53420.19
$ws.Range("K40").Value = 53420.19
$ws.Range("M40").Value = -53284.19
$ws.Range("H46").Value = 3304.3462
$ws.Range("I46").Value = 2332.1177
$ws.Range("K46").Value = 2332.1177
$ws.Range("M46").Value = -2144.1177
$ws.Range("H126").Value = 4869.5386
$ws.Range("I126").Value = 3855.375
$ws.Range("J126").Value = 6492.2
$ws.Range("K126").Value = 11566.125
$ws.Range("L126").Value = 19476.6
$ws.Range("M126").Value = -9096.125
$ws.Range("N126").Value = -24416.6
$ws.Range("H132").Value = 9496.385
$ws.Range("I132").Value = 7195.9
$ws.Range("K132").Value = 21587.7
$ws.Range("M132").Value = -19057.7
$ws.Range("H136").Value = 4323.4165
$ws.Range("I136").Value = 2150.1667
$ws.Range("K136").Value = 6450.500100000001
$ws.Range("M136").Value = -3900.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 29569.5
$ws.Range("I51").Value = 29569.5
$ws.Range("K51").Value = 29569.5
$ws.Range("M51").Value = -29059.5
$ws.Range("H53").Value = 30076
$ws.Range("I53").Value = 30076
$ws.Range("K53").Value = 30076
$ws.Range("M53").Value = -29469
$ws.Range("H122").Value = 19233730
$ws.Range("I122").Value = 3374.875
$ws.Range("J122").Value = 50002300
$ws.Range("K122").Value = 10124.625
$ws.Range("L122").Value = 150006900
$ws.Range("M122").Value = -7674.625
$ws.Range("N122").Value = -150011800
$ws.Range("H132").Value = 2734.311
$ws.Range("I132").Value = 2633.5898
$ws.Range("K132").Value = 7900.769400000001
$ws.Range("M132").Value = -5370.769400000001
$ws.Range("H136").Value = 3057.4312
$ws.Range("I136").Value = 2409.6047
$ws.Range("K136").Value = 7228.8141
$ws.Range("M136").Value = -4678.8141
